$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 296 (existing rows 296:314 shift down to 300:318).
$ws.Rows("296:299").Insert()

# --- New row 296 (Cuarta) ---
$ws.Range("A296").Value = 4
$ws.Range("B296").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C296").Value = "Los Lagos"
$ws.Range("D296").Value = 44931
$ws.Range("E296").Value = 10
$ws.Range("F296").Value = 100112028
$ws.Range("G296").Value = "Sandia"
$ws.Range("H296").Value = "Sin especificar"
$ws.Range("I296").Value = "Cuarta"
$ws.Range("J296").Value = 1000
$ws.Range("K296").Value = 1500
$ws.Range("L296").Value = 1500
$ws.Range("M296").Value = 1500
$ws.Range("N296").Value = "$/unidad"
$ws.Range("O296").Value = "Región del Maule"
$ws.Range("P296").Value = 1500
$ws.Range("Q296").Value = 1
$ws.Range("R296").Value = "Hortaliza"

# --- New row 297 (Primera) ---
$ws.Range("A297").Value = 4
$ws.Range("B297").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C297").Value = "Los Lagos"
$ws.Range("D297").Value = 44931
$ws.Range("E297").Value = 10
$ws.Range("F297").Value = 100112028
$ws.Range("G297").Value = "Sandia"
$ws.Range("H297").Value = "Sin especificar"
$ws.Range("I297").Value = "Primera"
$ws.Range("J297").Value = 1000
$ws.Range("K297").Value = 4000
$ws.Range("L297").Value = 4000
$ws.Range("M297").Value = 4000
$ws.Range("N297").Value = "$/unidad"
$ws.Range("O297").Value = "Región del Maule"
$ws.Range("P297").Value = 4000
$ws.Range("Q297").Value = 1
$ws.Range("R297").Value = "Hortaliza"

# --- New row 298 (Segunda) ---
$ws.Range("A298").Value = 4
$ws.Range("B298").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C298").Value = "Los Lagos"
$ws.Range("D298").Value = 44931
$ws.Range("E298").Value = 10
$ws.Range("F298").Value = 100112028
$ws.Range("G298").Value = "Sandia"
$ws.Range("H298").Value = "Sin especificar"
$ws.Range("I298").Value = "Segunda"
$ws.Range("J298").Value = 1000
$ws.Range("K298").Value = 3500
$ws.Range("L298").Value = 3500
$ws.Range("M298").Value = 3500
$ws.Range("N298").Value = "$/unidad"
$ws.Range("O298").Value = "Región del Maule"
$ws.Range("P298").Value = 3500
$ws.Range("Q298").Value = 1
$ws.Range("R298").Value = "Hortaliza"

# --- New row 299 (Tercera) ---
$ws.Range("A299").Value = 4
$ws.Range("B299").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C299").Value = "Los Lagos"
$ws.Range("D299").Value = 44931
$ws.Range("E299").Value = 10
$ws.Range("F299").Value = 100112028
$ws.Range("G299").Value = "Sandia"
$ws.Range("H299").Value = "Sin especificar"
$ws.Range("I299").Value = "Tercera"
$ws.Range("J299").Value = 1000
$ws.Range("K299").Value = 2500
$ws.Range("L299").Value = 2500
$ws.Range("M299").Value = 2500
$ws.Range("N299").Value = "$/unidad"
$ws.Range("O299").Value = "Región del Maule"
$ws.Range("P299").Value = 2500
$ws.Range("Q299").Value = 1
$ws.Range("R299").Value = "Hortaliza"
